$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two more quarterly periods (Dec-2018 and Sep-2018) to the
# --- CBM quarterly financials tables. The sheet holds three stacked
# --- tables (Income Statement, Balance Sheet, Cash Flow Statement),
# --- each with a "Period Ending" row followed by line items in
# --- columns D..K (8 quarters). We insert two new columns at D:E so
# --- the newest two quarters lead, pushing the existing 8 quarters
# --- of data right into F:M.
$ws.Columns("D:E").Insert()

# Carry the existing number formatting (date format for the header
# rows, plain number format for the data rows) from the columns that
# just got shifted right (F:G, which held the old D:E data/styles)
# into the freshly inserted, blank D:E columns.
$ws.Range("F7:G102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the column width Excel applied to the two new columns.
$ws.Columns("D:E").ColumnWidth = 13.830729166666666

# --- New quarter data (period endings 2018-12-31 and 2018-09-30)
# --- for all three tables, plus small corrected figures that moved
# --- along with the rest of the historical data.
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 134300
$ws.Range("E8").Value = 104600
$ws.Range("D9").Value = 86000
$ws.Range("E9").Value = 71900
$ws.Range("D10").Value = 48300
$ws.Range("E10").Value = 32700
$ws.Range("D12").Value = 3600
$ws.Range("E12").Value = 4200
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 3400
$ws.Range("E14").Value = 7400
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 114500
$ws.Range("E17").Value = 98000
$ws.Range("D18").Value = 19800
$ws.Range("E18").Value = 6600
$ws.Range("D20").Value = -1000
$ws.Range("E20").Value = 4800
$ws.Range("D21").Value = 32400
$ws.Range("E21").Value = 20600
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 18800
$ws.Range("E23").Value = 11400
$ws.Range("D24").Value = 17300
$ws.Range("E24").Value = -13300
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 1500
$ws.Range("E26").Value = 24700
$ws.Range("D27").Value = 1500
$ws.Range("E27").Value = 24700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -200
$ws.Range("E29").Value = 2000
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 1000
$ws.Range("E32").Value = -4800
$ws.Range("D33").Value = 1200
$ws.Range("E33").Value = 26700
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 1200
$ws.Range("E35").Value = 26700
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 95900
$ws.Range("E41").Value = 97100
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 185000
$ws.Range("E43").Value = 199300
$ws.Range("D44").Value = 111100
$ws.Range("E44").Value = 103600
$ws.Range("D45").Value = 18200
$ws.Range("E45").Value = 16800
$ws.Range("D46").Value = 410100
$ws.Range("E46").Value = 416800
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 360500
$ws.Range("E48").Value = 352900
$ws.Range("D49").Value = 448300
$ws.Range("E49").Value = 463400
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 4500
$ws.Range("E52").Value = 14700
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1223400
$ws.Range("E54").Value = 1247900
$ws.Range("D57").Value = 47000
$ws.Range("E57").Value = 38200
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 57400
$ws.Range("E59").Value = 53600
$ws.Range("D60").Value = 104400
$ws.Range("E60").Value = 91800
$ws.Range("D61").Value = 300000
$ws.Range("E61").Value = 325000
$ws.Range("D62").Value = 165300
$ws.Range("E62").Value = 171800
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 569700
$ws.Range("E66").Value = 588600
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 538500
$ws.Range("E72").Value = 537300
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 653700
$ws.Range("E76").Value = 659300
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 1200
$ws.Range("E81").Value = 26700
$ws.Range("D83").Value = 13700
$ws.Range("E83").Value = 9200
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 47200
$ws.Range("E89").Value = 22300
$ws.Range("D91").Value = -19100
$ws.Range("E91").Value = -10800
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -24300
$ws.Range("E94").Value = -429800
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -23600
$ws.Range("E100").Value = 333100
$ws.Range("D101").Value = -500
$ws.Range("E101").Value = 200
$ws.Range("D102").Value = -1300
$ws.Range("E102").Value = -74200
# --- A handful of historical figures (now in F:I) were corrected
# --- at the same time as the new-quarter data was added.
$ws.Range("G14").Value = 300
$ws.Range("I14").Value = 0
$ws.Range("H17").Value = 126000
$ws.Range("I17").Value = 87100
$ws.Range("H18").Value = 56300
$ws.Range("I18").Value = 25500
$ws.Range("H20").Value = -500
$ws.Range("I20").Value = -800
$ws.Range("H32").Value = 500
$ws.Range("I32").Value = 800